$wb = $excel.ActiveWorkbook

# --- Recreate the "05 - Add Employee" sheet so it gets a fresh sheetId ---
# (mirrors the diff: sheetId bumps 4 -> 5 while keeping the same name/position/rId)
$oldWs = $wb.Worksheets.Item("05 - Add Employee")
$lastWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add($null, $lastWs)
$newIndex = $newWs.Index
$oldWs.Delete()
$ws = $wb.Worksheets.Item($newIndex - 1)
$ws.Name = "05 - Add Employee"

# --- Row 1: headers ---
$ws.Range("A1").Value = "empFirstName"
$ws.Range("B1").Value = "empMidName"
$ws.Range("C1").Value = "empLastName"
$ws.Range("D1").Value = "empID"
$ws.Range("E1").Value = "empUsername"
$ws.Range("F1").Value = "empPassword"
$ws.Range("G1").Value = "empPasswordConfirmaton"
$ws.Range("H1").Value = "status"
$ws.Range("I1").Value = "file"

# --- Row 2: Ahmad Tim Sean ---
$ws.Range("A2").Value = "Ahmad"
$ws.Range("B2").Value = "Tim"
$ws.Range("C2").Value = "Sean"
$ws.Range("D2").Value = "'0299"
$ws.Range("E2").Value = "ahmad098"
$ws.Range("F2").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("G2").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("H2").Value = 1
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("I2").Value = "C:\\Users\\Hasnul\\Katalon Studio\\Katalon-Studio-Training-Assessment\\OrangeHRM_Nafis\\File Upload Test Data\\19263862.png"

# --- Row 3: Rick Froberg ---
$ws.Range("A3").Value = "Rick"
$ws.Range("C3").Value = "Froberg"
$ws.Range("D3").Value = "'0322"
$ws.Range("E3").Value = "rickforob166"
$ws.Range("F3").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("G3").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("H3").Value = 2
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("I3").Value = "C:\\Users\\Hasnul\\Katalon Studio\\Katalon-Studio-Training-Assessment\\OrangeHRM_Nafis\\File Upload Test Data\\43392873.png"

# --- Row 4: Alfred Bok Bok ---
$ws.Range("A4").Value = "Alfred"
$ws.Range("C4").Value = "Bok Bok"
$ws.Range("D4").Value = "'0666"
$ws.Range("E4").Value = "alfredbok2"
$ws.Range("F4").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("G4").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "C:\\Users\\Hasnul\\Katalon Studio\\Katalon-Studio-Training-Assessment\\OrangeHRM_Nafis\\File Upload Test Data\\19263862.png"

# --- Row 5: Zack Junaidi Apron ---
$ws.Range("A5").Value = "Zack"
$ws.Range("B5").Value = "Junaidi"
$ws.Range("C5").Value = "Apron"
$ws.Range("D5").Value = "'1066"
$ws.Range("E5").Value = "zackapron99"
$ws.Range("F5").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("G5").Value = "jq6leojyGnb59OuxS61Hr0q+yqPfNELT"
$ws.Range("H5").Value = 2

$ws.Range("E9").Select()
